# Update imputed values produced by the RandomForest algorithm run.
# (commit message: "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 9.020199999999996
$ws.Range("C4").Value = -11.98219999999999
$ws.Range("B6").Value = 5.429100000000002
$ws.Range("B7").Value = 5.344899999999999
$ws.Range("C9").Value = -10.2262
$ws.Range("C12").Value = -10.58169999999999
$ws.Range("B16").Value = 4.9396
$ws.Range("C17").Value = -14.39609999999998
$ws.Range("C18").Value = -11.94959999999999
$ws.Range("C19").Value = -11.2685
$ws.Range("B20").Value = 9.657999999999989
$ws.Range("C20").Value = -12.5254
$ws.Range("C26").Value = -11.9207
$ws.Range("B28").Value = 5.829600000000004
$ws.Range("B29").Value = 5.445500000000004
$ws.Range("C31").Value = -13.4529
$ws.Range("B32").Value = 6.700399999999998
$ws.Range("C39").Value = -11.7525
$ws.Range("B40").Value = 9.176499999999988
$ws.Range("C40").Value = -12.66670000000001
$ws.Range("C41").Value = -12.43650000000001
$ws.Range("C42").Value = -11.922
$ws.Range("C43").Value = -13.16759999999999
$ws.Range("B46").Value = 5.413100000000003
$ws.Range("C47").Value = -12.18159999999999
$ws.Range("C48").Value = -12.169
$ws.Range("B51").Value = 5.399699999999998
$ws.Range("B52").Value = 5.436399999999996
$ws.Range("B57").Value = 5.420499999999996
$ws.Range("B59").Value = 4.502800000000002
$ws.Range("B62").Value = 5.144499999999996
$ws.Range("C63").Value = -10.1649
$ws.Range("C64").Value = -10.38819999999999
$ws.Range("B66").Value = 5.197700000000001
$ws.Range("B73").Value = 8.810399999999996
$ws.Range("B74").Value = 8.971099999999993
$ws.Range("C76").Value = -12.2085
$ws.Range("C81").Value = -14.18109999999999
$ws.Range("C89").Value = -13.6525
$ws.Range("B92").Value = 4.718999999999999
$ws.Range("C94").Value = -10.6598
$ws.Range("B100").Value = 4.732300000000003
